# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45205 (2023-10-06) to 45206 (2023-10-07).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 452
$firstRow = 2

$range = $ws.Range("C$firstRow`:C$lastRow")
$range.Value = 45206
